$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ 'C'=13.37041050194737; 'D'=4.155204696756664; 'E'=13.72610928902245; 'F'=26.5227943986173; 'G'=34.09956198786418; 'H'=14.5584433102773; 'L'=9.145448808988293; 'N'=18.99769123883711; 'O'=23.13022551420616 }
    3 = @{ 'C'=13.26998484690672; 'D'=4.164241322338348; 'E'=13.66501597118269; 'F'=26.12679731536193; 'G'=33.26842763544568; 'H'=14.50374877142435; 'L'=9.149145858536359; 'N'=18.40031508502701; 'O'=22.89358244757704 }
    4 = @{ 'C'=13.21156305405411; 'D'=4.169974493241716; 'E'=13.63080591497906; 'F'=25.88873046564321; 'G'=32.75873017636872; 'H'=14.47347989785808; 'L'=9.153051092502629; 'N'=18.02485520896362; 'O'=22.75360380579077 }
    5 = @{ 'C'=13.18859190932568; 'D'=4.17235756564671; 'E'=13.61770494133458; 'F'=25.79312095759196; 'G'=32.55152567640876; 'H'=14.46198665716708; 'L'=9.155053650312562; 'N'=17.86990355188765; 'O'=22.69796028144011 }
    6 = @{ 'C'=13.18482866954959; 'D'=4.172756106705926; 'E'=13.61558052323208; 'F'=25.77733359755572; 'G'=32.51716006840205; 'H'=14.46012926374611; 'L'=9.155411001309862; 'N'=17.84406337566597; 'O'=22.68880692487065 }
    7 = @{ 'C'=13.21124984418041; 'D'=4.170006442449698; 'E'=13.63062581798043; 'F'=25.88743518790102; 'G'=32.75593325152816; 'H'=14.47332147840509; 'L'=9.153076435169154; 'N'=18.02277304767603; 'O'=22.75284763561852 }
    8 = @{ 'C'=13.33512466194454; 'D'=4.158282400405042; 'E'=13.70436527329524; 'F'=26.38527843459258; 'G'=33.81306139919632; 'H'=14.53890214088213; 'L'=9.146384169886828; 'N'=18.79364780656867; 'O'=23.04756400656899 }
    9 = @{ 'C'=13.60273188941015; 'D'=4.136741472072262; 'E'=13.87468339608541; 'F'=27.39610705352802; 'G'=35.87645284445449; 'H'=14.69338825875396; 'L'=9.146232988110823; 'N'=20.2273683202997; 'O'=23.66491092064396 }
    10 = @{ 'C'=13.81293516011308; 'D'=4.121777607462343; 'E'=14.01477531325477; 'F'=28.15212020161805; 'G'=37.36806833338232; 'H'=14.82205274249375; 'L'=9.154019931566735; 'N'=21.22223697909767; 'O'=24.13857679840676 }
    11 = @{ 'C'=13.91117940161917; 'D'=4.115152714106183; 'E'=14.08158059664049; 'F'=28.49733257179048; 'G'=38.03778817837413; 'H'=14.88372749745571; 'L'=9.159272117890167; 'N'=21.66018057919901; 'O'=24.3575587741136 }
    12 = @{ 'C'=13.94872960505786; 'D'=4.112669886466443; 'E'=14.10730431498528; 'F'=28.62811052037343; 'G'=38.28985415784105; 'H'=14.90752024431413; 'L'=9.161506106580909; 'N'=21.82377585682186; 'O'=24.44091257073001 }
    13 = @{ 'C'=13.94062753326048; 'D'=4.113203462258511; 'E'=14.10174557542708; 'F'=28.5999448485599; 'G'=38.2356399765688; 'H'=14.90237681339962; 'L'=9.161014090389759; 'N'=21.78864458690801; 'O'=24.42294289741484 }
    14 = @{ 'C'=13.91426189276528; 'D'=4.114947933691749; 'E'=14.08368846160805; 'F'=28.50809142899396; 'G'=38.05855820458893; 'H'=14.88567625069102; 'L'=9.159451002355464; 'N'=21.67368539489659; 'O'=24.3644081484729 }
    15 = @{ 'C'=13.89815650436399; 'D'=4.11601983276915; 'E'=14.07268293517114; 'F'=28.45183163900173; 'G'=37.94988186959795; 'H'=14.87550325658976; 'L'=9.158525459154946; 'N'=21.60297336126124; 'O'=24.32860772178101 }
    16 = @{ 'C'=13.80656489398851; 'D'=4.12221419876783; 'E'=14.01046995804423; 'F'=28.12957441757639; 'G'=37.32409916418752; 'H'=14.8180842227151; 'L'=9.153711033056302; 'N'=21.19330956972086; 'O'=24.12433005285303 }
    17 = @{ 'C'=13.75102692519575; 'D'=4.126060683422856; 'E'=13.97308092542864; 'F'=27.93211667036676; 'G'=36.9377407366258; 'H'=14.78365497450506; 'L'=9.151194975731659; 'N'=20.93814219015166; 'O'=23.99985807693567 }
    18 = @{ 'C'=13.71933155048363; 'D'=4.128290253099933; 'E'=13.95186638201738; 'F'=27.81867415254562; 'G'=36.71470216004653; 'H'=14.7641490806374; 'L'=9.149908803807708; 'N'=20.79000725568362; 'O'=23.92859985185283 }
    19 = @{ 'C'=13.70864363409617; 'D'=4.129048106701211; 'E'=13.94473390647452; 'F'=27.78029090189914; 'G'=36.63905374850814; 'H'=14.75759613660141; 'L'=9.149501000575199; 'N'=20.73962067985785; 'O'=23.90453280603581 }
    20 = @{ 'C'=13.7569135216573; 'D'=4.125649443604662; 'E'=13.97703108809821; 'F'=27.95312391127031; 'G'=36.97895569104387; 'H'=14.78728940031413; 'L'=9.151446156807499; 'N'=20.96544799484619; 'O'=24.01307422436261 }
    21 = @{ 'C'=13.92199693653972; 'D'=4.114434840522379; 'E'=14.08898084690295; 'F'=28.53507059482261; 'G'=38.11061535280545; 'H'=14.89056984565117; 'L'=9.159903474271863; 'N'=21.70751365554066; 'O'=24.38159014469124 }
    22 = @{ 'C'=14.03189830191663; 'D'=4.107256118724323; 'E'=14.16462140604692; 'F'=28.91565940356294; 'G'=38.84111764272446; 'H'=14.96061455826228; 'L'=9.166858972503363; 'N'=22.17935961385674; 'O'=24.62490897187242 }
    23 = @{ 'C'=13.97306795110939; 'D'=4.111073860734327; 'E'=14.12402979452875; 'F'=28.7125518170347; 'G'=38.45215282951914; 'H'=14.92300244718488; 'L'=9.163016315270088; 'N'=21.92877110911181; 'O'=24.49484349743257 }
    24 = @{ 'C'=13.75425145984138; 'D'=4.125835308476042; 'E'=13.97524434334382; 'F'=27.94362628268672; 'G'=36.96032524213635; 'H'=14.78564537839181; 'L'=9.151332098343742; 'N'=20.95310750188673; 'O'=24.00709825654724 }
    25 = @{ 'C'=13.52783970910533; 'D'=4.142415879446898; 'E'=13.82592466555299; 'F'=27.11974928405018; 'G'=35.32125530339712; 'H'=14.64888441999562; 'L'=9.144885548712821; 'N'=19.84905939529497; 'O'=23.13022551420616 }
}

foreach ($r in $data.Keys) {
    $rowData = $data[$r]
    foreach ($col in $rowData.Keys) {
        $ws.Range("$col$r").Value = $rowData[$col]
    }
}
